$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CF (column B) values for rows 2-94, derived from the target diff.
$values = @(
    0.02796140685677528,
    0.02796140685677528,
    0.02796140685677528,
    0.02796140685677528,
    0.04261558502912521,
    0.03100503422319889,
    0.07880202680826187,
    0.1343789994716644,
    0.1469141393899918,
    0.1000665202736855,
    0.03529689833521843,
    0.02820187248289585,
    0.04917032644152641,
    0.05099577829241753,
    0.04806243255734444,
    0.1717489510774612,
    0.3450882434844971,
    0.3895267844200134,
    0.3939395546913147,
    0.3837548196315765,
    0.3795405924320221,
    0.3858906328678131,
    0.3736143708229065,
    0.3833216726779938,
    0.3875951766967773,
    0.3941161930561066,
    0.3946114480495453,
    0.394012063741684,
    0.3729609847068787,
    0.3054832220077515,
    0.2029432356357574,
    0.08896761387586594,
    0.05902456119656563,
    0.0827145054936409,
    0.07151544839143753,
    0.1140044629573822,
    0.1185560449957848,
    0.1265325099229813,
    0.1185560449957848,
    0.07957780361175537,
    0.04900955781340599,
    0.04671235382556915,
    0.04671235382556915,
    0.04671235382556915,
    0.04671235382556915,
    0.0408606268465519,
    0.03269428387284279,
    0.02796140685677528,
    0.03549665212631226,
    0.04261558502912521,
    0.03549665212631226,
    0.02796140685677528,
    0.02796140685677528,
    0.02856598980724812,
    0.02856598980724812,
    0.02796140685677528,
    0.02796140685677528,
    0.02820187248289585,
    0.02820187248289585,
    0.02820187248289585,
    0.02820187248289585,
    0.02820187248289585,
    0.02820187248289585,
    0.02820187248289585,
    0.02796140685677528,
    0.02796140685677528,
    0.03539891168475151,
    0.08561854809522629,
    0.1005492806434631,
    0.09617377817630768,
    0.0614122599363327,
    0.04922284558415413,
    0.04117486998438835,
    0.04647859930992126,
    0.04647859930992126,
    0.07439888268709183,
    0.0592782236635685,
    0.04963897541165352,
    0.03819409385323524,
    0.03543885797262192,
    0.0371975488960743,
    0.03068605065345764,
    0.03068605065345764,
    0.03548388183116913,
    0.03743801638484001,
    0.03743801638484001,
    0.03264018148183823,
    0.03034817613661289,
    0.02796140685677528,
    0.02796140685677528,
    0.02796140685677528,
    0.02796140685677528,
    0.02796140685677528
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
